# Insert a new price-record row at row 214 (pushing the existing rows
# 214-248 down to 215-249) and populate it with the new weekly data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 214..248 down to 215..249, carrying formatting (e.g. the date
# style on column D) along with them.
$ws.Rows.Item(214).Insert()

# Populate the newly inserted row 214 with the new record.
$ws.Range("A214").Value = 6
$ws.Range("B214").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C214").Value = "Metropolitana"
$ws.Range("D214").Value = 44505
$ws.Range("E214").Value = 13
$ws.Range("F214").Value = 100112032
$ws.Range("G214").Value = "Zapallo italiano"
$ws.Range("H214").Value = "Sin especificar"
$ws.Range("I214").Value = "Primera"
$ws.Range("J214").Value = 400
$ws.Range("K214").Value = 6000
$ws.Range("L214").Value = 7000
$ws.Range("M214").Value = 6425
$ws.Range("N214").Value = "`$/caja 50 unidades"
$ws.Range("O214").Value = "Región de O'Higgins"
$ws.Range("P214").Value = 128
$ws.Range("Q214").Value = 50
$ws.Range("R214").Value = "Hortaliza"
